$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteSpecial constants used below to duplicate a source cell's
# value and formatting onto a new cell without inventing new style
# combinations (property-by-property copies like .Interior.Color /
# .Font.Bold tend to create brand-new font/fill entries in styles.xml).
$xlPasteValues  = -4163
$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1. New header pair: AP1 "03-02_A" / AQ1 "03-02_0" - clones of the
#    existing AN1 / AO1 headers (same bold/centered header style),
#    just with the next day's label.
# ---------------------------------------------------------------------
$ws.Range("AN1").Copy()
$ws.Range("AP1").PasteSpecial($xlPasteFormats)
$ws.Range("AP1").Value2 = "03-02_A"

$ws.Range("AO1").Copy()
$ws.Range("AQ1").PasteSpecial($xlPasteFormats)
$ws.Range("AQ1").Value2 = "03-02_0"

# ---------------------------------------------------------------------
# 2. Every data row: duplicate the "03-01" pair (AN/AO) into the new
#    "03-02" pair (AP/AQ), and normalize the old AO column from a
#    text-that-looks-like-a-number into a real number (AQ keeps the
#    text form the way AO used to have it).
# ---------------------------------------------------------------------
$lastRow = 173
for ($r = 2; $r -le $lastRow; $r++) {
    $an = $ws.Range("AN$r")
    $ao = $ws.Range("AO$r")
    $ap = $ws.Range("AP$r")
    $aq = $ws.Range("AQ$r")

    # AQ <- clone of the current (still text) AO, before AO is touched.
    $ao.Copy()
    $aq.PasteSpecial($xlPasteValues)
    $ao.Copy()
    $aq.PasteSpecial($xlPasteFormats)

    # AP <- clone of AN (value + fill/border/font style).
    $an.Copy()
    $ap.PasteSpecial($xlPasteValues)
    $an.Copy()
    $ap.PasteSpecial($xlPasteFormats)

    # AO, in place: re-assigning Value2 lets the host re-detect a
    # numeric-looking string as a genuine number. Skip truly blank
    # rows - round-tripping an empty value would delete the cell.
    $aoVal = $ao.Value2
    if ($aoVal -ne $null -and $aoVal -ne "") {
        $ao.Value2 = $aoVal
    }
}

# ---------------------------------------------------------------------
# 3. Column A for rows 171-173 had been stored as text ("41231396", …)
#    even though every other row already held a real number; the same
#    Value2 round-trip normalizes them in place.
# ---------------------------------------------------------------------
foreach ($r in 171..173) {
    $a = $ws.Range("A$r")
    $a.Value2 = $a.Value2
}
